$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = 'questions = [
    {
        "title": "Imagine you are developing a word processor that represents text as linked lists. In this system, every character in a text string is represented by a node in the linked list. For example, the word \"Hello\" is represented as a linked list: ''H'' -&gt ''e'' -&gt ''l'' -&gt ''l'' -&gt ''o'' -&gt Null.Your task is to implement a function string_list_slice that allows users to select a portion of the text, starting at a specified index and ending at another specified index. So, if we specify 0 and 4, our slice will consist of the values at indices 0, 1, 2, and 3. This selected portion should be returned as a new linked list, representing the selected text. The linked list stores the string like so:The inputs to the string_list_slice function will be string_list (a string), begin (the starting index), and end (the ending index + 1). The function should return a string containing the requested slice (substring).Note:The index values supplied assume that the programming language is 0-indexed. If you choose to solve the problem in a language that is 1-indexed, you must take that into account and make the adjustment.If the ending index exceeds the end of the list, include characters up to the end.The function should work for linked lists representing strings of any length.Example 1Input:string_list = StringList(''A whole new internet'')begin = 2end = 11Output:''whole new''Explanation:We begin at index 2 (the first letter \u2018w\u2019) and continue until index 10 (the second letter ''w'')Example 2Input:string_list = StringList(''Computing calmly at home'')begin = 0end = 50Output:''Computing calmly at home''Explanation:Because our ending index exceeds the length of the string, we just stop at the end",
        "ques_type": null,
        "options": [],
        "score": null
    }
]'
$ws.Rows(2).Delete()
$ws.Rows(1).AutoFit()
